# Update gh-pages to output generated at efdf629
#
# The "最低票价" (lowest ticket price) column (G) on the "展览" and
# "全部类型" sheets moves from raw numeric yuan-fen values to the
# human-readable strings actually shown on the source site (prices in
# whole yuan, or status text like "已售罄"/"不可售"/"预售中" when there is
# no numeric price). A handful of "想去人数" (interest count) values in
# column F were also refreshed to newer counts.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$cellRef, [string]$text)
    $c = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the value as text even when
    # it looks numeric (e.g. "60"); re-applying the Normal style afterwards
    # drops the implicit "quote prefix" formatting so only the cell's value
    # (not its style) changes.
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

# (row, new G text, optional new F number)
$rows = @(
    @{ Row = 2;  G = "60";    F = $null },
    @{ Row = 3;  G = "不可售"; F = $null },
    @{ Row = 4;  G = "58";    F = $null },
    @{ Row = 5;  G = "25";    F = 10455 },
    @{ Row = 7;  G = "已售罄"; F = $null },
    @{ Row = 8;  G = "55";    F = $null },
    @{ Row = 9;  G = "50";    F = $null },
    @{ Row = 10; G = "49";    F = $null },
    @{ Row = 11; G = "45";    F = 13 },
    @{ Row = 12; G = "60";    F = $null },
    @{ Row = 13; G = "68";    F = $null },
    @{ Row = 14; G = "65";    F = $null },
    @{ Row = 15; G = "78";    F = $null },
    @{ Row = 16; G = "48";    F = $null },
    @{ Row = 17; G = "预售中"; F = 10605 },
    @{ Row = 18; G = "65";    F = $null },
    @{ Row = 19; G = "1";     F = $null },
    @{ Row = 20; G = "1";     F = $null },
    @{ Row = 21; G = "1";     F = $null },
    @{ Row = 22; G = "1";     F = $null },
    @{ Row = 23; G = "1";     F = $null },
    @{ Row = 24; G = "1";     F = $null }
)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in $rows) {
        $gref = "G$($r.Row)"
        Set-TextValue $ws $gref $r.G

        if ($null -ne $r.F) {
            $fref = "F$($r.Row)"
            $ws.Range($fref).Value = $r.F
        }
    }
}

# Row 6's "想去人数" count diverged between the two sheets when the data
# was refreshed.
$wsExhibit = $wb.Worksheets.Item("展览")
Set-TextValue $wsExhibit "G6" "60"
$wsExhibit.Range("F6").Value = 9072

$wsAllTypes = $wb.Worksheets.Item("全部类型")
Set-TextValue $wsAllTypes "G6" "60"
$wsAllTypes.Range("F6").Value = 9073

Write-Host "Applied updates to Exhibitions and AllTypes sheets"
